$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-2.787195682525635, 8.478515625, 9.063291549682615),
    @(-10.50285053253174, -24.10493659973145, -11.3930606842041),
    @(-2.336788654327393, 8.028789520263672, -11.69052314758301),
    @(-33.8237419128418, -34.27063751220703, -13.57823276519775),
    @(-2.554568767547607, -39.26761627197266, 14.68596649169922),
    @(-6.776285171508789, 9.06696891784668, -3.473462343215942),
    @(-4.451385021209717, -37.8787841796875, 33.21585845947266),
    @(-25.62279510498047, 7.356860160827637, -15.05502891540527),
    @(-8.065028190612793, -2.866869926452637, -9.438656806945801),
    @(26.9844970703125, -14.56967449188232, 10.2812442779541)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
